$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3099
$ws.Range("J70").Value = 3218.8
$ws.Range("L70").Value = 9656.400000000001
$ws.Range("N70").Value = -10196.4
$ws.Range("H73").Value = 3099
$ws.Range("J73").Value = 3218.8
$ws.Range("L73").Value = 9656.400000000001
$ws.Range("N73").Value = -11528.4
$ws.Range("H86").Value = 934.75
$ws.Range("I86").Value = 854
$ws.Range("K86").Value = 854
$ws.Range("M86").Value = 269
$ws.Range("H89").Value = 934.75
$ws.Range("I89").Value = 854
$ws.Range("K89").Value = 4270
$ws.Range("M89").Value = 1346
$ws.Range("H116").Value = 4224.933
$ws.Range("I116").Value = 3952.6667
$ws.Range("J116").Value = 4633.3335
$ws.Range("K116").Value = 3952.6667
$ws.Range("L116").Value = 4633.3335
$ws.Range("M116").Value = -510.6667000000002
$ws.Range("N116").Value = -11517.3335
$ws.Range("H138").Value = 3948.457
$ws.Range("J138").Value = 5006.7144
$ws.Range("L138").Value = 15020.1432
$ws.Range("N138").Value = -25300.1432
$ws.Range("H141").Value = 2602.3
$ws.Range("I141").Value = 2669.3333
$ws.Range("K141").Value = 8007.999899999999
$ws.Range("M141").Value = -2827.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3134.682
$ws.Range("I45").Value = 2530.9167
$ws.Range("J45").Value = 3859.2
$ws.Range("K45").Value = 2530.9167
$ws.Range("L45").Value = 3859.2
$ws.Range("M45").Value = -2153.9167
$ws.Range("N45").Value = -4613.2
$ws.Range("H46").Value = 8112.8184
$ws.Range("I46").Value = 6298
$ws.Range("K46").Value = 6298
$ws.Range("M46").Value = -5979
$ws.Range("H61").Value = 2639.5938
$ws.Range("I61").Value = 2615.5667
$ws.Range("K61").Value = 2615.5667
$ws.Range("M61").Value = -2403.5667
$ws.Range("H63").Value = 1430575
$ws.Range("I63").Value = 2198
$ws.Range("J63").Value = 5001517.5
$ws.Range("K63").Value = 2198
$ws.Range("L63").Value = 5001517.5
$ws.Range("M63").Value = -1512
$ws.Range("N63").Value = -5002889.5
$ws.Range("H66").Value = 1430575
$ws.Range("I66").Value = 2198
$ws.Range("J66").Value = 5001517.5
$ws.Range("K66").Value = 10990
$ws.Range("L66").Value = 25007587.5
$ws.Range("M66").Value = -7558
$ws.Range("N66").Value = -25014451.5
$ws.Range("H74").Value = 1428.9333
$ws.Range("I74").Value = 1388.1428
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1388.1428
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -514.1428000000001
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 1428.9333
$ws.Range("I77").Value = 1388.1428
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 6940.714
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -2572.714
$ws.Range("N77").Value = -18736
$ws.Range("H102").Value = 2611.5
$ws.Range("I102").Value = 1778.8
$ws.Range("J102").Value = 6775
$ws.Range("K102").Value = 1778.8
$ws.Range("L102").Value = 6775
$ws.Range("M102").Value = -156.8
$ws.Range("N102").Value = -10019
$ws.Range("H131").Value = 52452.4
$ws.Range("J131").Value = 52452.4
$ws.Range("L131").Value = 52452.4
$ws.Range("N131").Value = -62532.4
$ws.Range("H132").Value = 1418.138
$ws.Range("I132").Value = 1478
$ws.Range("J132").Value = 899.3333
$ws.Range("K132").Value = 4434
$ws.Range("L132").Value = 2697.9999
$ws.Range("M132").Value = -1904
$ws.Range("N132").Value = -7757.9999
$ws.Range("H135").Value = 94238.836
$ws.Range("J135").Value = 94238.836
$ws.Range("L135").Value = 94238.836
$ws.Range("N135").Value = -104378.836
$ws.Range("H136").Value = 2639.5938
$ws.Range("I136").Value = 2615.5667
$ws.Range("K136").Value = 7846.7001
$ws.Range("M136").Value = -5296.7001
$ws.Range("H141").Value = 84714.5
$ws.Range("J141").Value = 84714.5
$ws.Range("L141").Value = 84714.5
$ws.Range("N141").Value = -95074.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 23351.6
$ws.Range("I96").Value = 23351.6
$ws.Range("K96").Value = 23351.6
$ws.Range("M96").Value = -20605.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32261248
$ws.Range("I31").Value = 47620696
$ws.Range("K31").Value = 47620696
$ws.Range("M31").Value = -47620401
$ws.Range("H34").Value = 32261248
$ws.Range("I34").Value = 47620696
$ws.Range("K34").Value = 47620696
$ws.Range("M34").Value = -47620494
$ws.Range("H86").Value = 7188.5264
$ws.Range("I86").Value = 7137.9375
$ws.Range("J86").Value = 7458.3335
$ws.Range("K86").Value = 7137.9375
$ws.Range("L86").Value = 7458.3335
$ws.Range("M86").Value = -6014.9375
$ws.Range("N86").Value = -9704.333500000001
$ws.Range("H89").Value = 7188.5264
$ws.Range("I89").Value = 7137.9375
$ws.Range("J89").Value = 7458.3335
$ws.Range("K89").Value = 35689.6875
$ws.Range("L89").Value = 37291.6675
$ws.Range("M89").Value = -30073.6875
$ws.Range("N89").Value = -48523.6675
$ws.Range("H94").Value = 1768.75
$ws.Range("I94").Value = 2195
$ws.Range("J94").Value = 1342.5
$ws.Range("K94").Value = 2195
$ws.Range("L94").Value = 1342.5
$ws.Range("M94").Value = -1744
$ws.Range("N94").Value = -2244.5
$ws.Range("H141").Value = 117650.95
$ws.Range("J141").Value = 126101.3
$ws.Range("L141").Value = 126101.3
$ws.Range("N141").Value = -136461.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 100549.8
$ws.Range("I34").Value = 624.75
$ws.Range("K34").Value = 1874.25
$ws.Range("M34").Value = -1790.25
$ws.Range("H107").Value = 715229.1
$ws.Range("J107").Value = 910074.4399999999
$ws.Range("L107").Value = 2730223.32
$ws.Range("N107").Value = -2734063.32
$ws.Range("H132").Value = 27778322
$ws.Range("J132").Value = 598.1539
$ws.Range("L132").Value = 5383.3851
$ws.Range("N132").Value = -10443.3851

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 51795
$ws.Range("J46").Value = 51795
$ws.Range("L46").Value = 51795
$ws.Range("N46").Value = -52107
$ws.Range("H80").Value = 2545.158
$ws.Range("I80").Value = 2384.1
$ws.Range("J80").Value = 2724.111
$ws.Range("K80").Value = 2384.1
$ws.Range("L80").Value = 2724.111
$ws.Range("M80").Value = -1386.1
$ws.Range("N80").Value = -4720.111
$ws.Range("H83").Value = 2545.158
$ws.Range("I83").Value = 2384.1
$ws.Range("J83").Value = 2724.111
$ws.Range("K83").Value = 11920.5
$ws.Range("L83").Value = 13620.555
$ws.Range("M83").Value = -6928.5
$ws.Range("N83").Value = -23604.555
$ws.Range("H122").Value = 2778.276
$ws.Range("I122").Value = 2660.4783
$ws.Range("J122").Value = 3229.8333
$ws.Range("K122").Value = 7981.4349
$ws.Range("L122").Value = 9689.499899999999
$ws.Range("M122").Value = -5531.4349
$ws.Range("N122").Value = -14589.4999
$ws.Range("H132").Value = 1840.5385
$ws.Range("I132").Value = 1008.5862
$ws.Range("J132").Value = 4253.2
$ws.Range("K132").Value = 3025.7586
$ws.Range("L132").Value = 12759.6
$ws.Range("M132").Value = -495.7586000000001
$ws.Range("N132").Value = -17819.6
$ws.Range("H136").Value = 53247.777
$ws.Range("J136").Value = 53247.777
$ws.Range("L136").Value = 159743.331
$ws.Range("N136").Value = -164843.331

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 42616.5
$ws.Range("I61").Value = 48532.332
$ws.Range("J61").Value = 1205.6666
$ws.Range("K61").Value = 48532.332
$ws.Range("L61").Value = 1205.6666
$ws.Range("M61").Value = -48330.332
$ws.Range("N61").Value = -1609.6666
$ws.Range("H113").Value = 42616.5
$ws.Range("I113").Value = 48532.332
$ws.Range("J113").Value = 1205.6666
$ws.Range("K113").Value = 48532.332
$ws.Range("L113").Value = 1205.6666
$ws.Range("M113").Value = -46362.332
$ws.Range("N113").Value = -5545.6666
$ws.Range("H132").Value = 10564.475
$ws.Range("I132").Value = 10294.056
$ws.Range("K132").Value = 30882.168
$ws.Range("M132").Value = -28352.168

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 41978
$ws.Range("J99").Value = 41978
$ws.Range("L99").Value = 41978
$ws.Range("N99").Value = -47968
$ws.Range("H110").Value = 93333.336
$ws.Range("J110").Value = 93333.336
$ws.Range("L110").Value = 93333.336
$ws.Range("N110").Value = -101513.336
$ws.Range("H113").Value = 1509.96
$ws.Range("I113").Value = 1149.8572
$ws.Range("K113").Value = 3449.5716
$ws.Range("M113").Value = -1279.5716
$ws.Range("H122").Value = 5293.125
$ws.Range("I122").Value = 5116.294
$ws.Range("J122").Value = 5722.5713
$ws.Range("K122").Value = 15348.882
$ws.Range("L122").Value = 17167.7139
$ws.Range("M122").Value = -12898.882
$ws.Range("N122").Value = -22067.7139
$ws.Range("H132").Value = 3853.0244
$ws.Range("I132").Value = 3627.5938
$ws.Range("J132").Value = 4654.5557
$ws.Range("K132").Value = 10882.7814
$ws.Range("L132").Value = 13963.6671
$ws.Range("M132").Value = -8352.7814
$ws.Range("N132").Value = -19023.6671
